$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.866.91'
$ws.Range("E2").Value = '  -0.97%  '
$ws.Range("D3").Value = '1.562.78'
$ws.Range("E3").Value = '  +0.03%  '
$ws.Range("E4").Value = '  -0.21%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '205.95'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.18%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.488'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -1.00%  '
$ws.Range("E7").Value = '  -0.19%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '21.69'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -1.73%  '
$ws.Range("E9").Value = '  -0.11%  '
$ws.Range("E10").Value = '  -1.34%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0865'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.49%  '
$ws.Range("E12").Value = '  +0.05%  '
$ws.Range("D13").Value = '1.582.54'
$ws.Range("E13").Value = '  +0.83%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.72'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -1.07%  '
$ws.Range("E15").Value = '  -0.27%  '
$ws.Range("D16").Value = '26.879.62'
$ws.Range("E16").Value = '  -1.01%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.24'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -2.86%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '214.58'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +1.13%  '
$ws.Range("E19").Value = '  +2.01%  '
$ws.Range("E20").Value = '  -1.28%  '
$ws.Range("E21").Value = '  -0.20%  '
$ws.Range("E22").Value = '  +0.47%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.17'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -2.41%  '
$ws.Range("E24").Value = '  +1.10%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '154.00'
$ws.Range("D25").ClearFormats()
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.74'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +1.88%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '14.92'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.49%  '
$ws.Range("E28").Value = '  -0.19%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0463'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -0.02%  '
$ws.Range("E31").Value = '  -3.06%  '
$ws.Range("E32").Value = '  -0.04%  '
$ws.Range("D33").Value = '1.402.15'
$ws.Range("E33").Value = '  +1.92%  '
$ws.Range("E34").Value = '  -0.66%  '
$ws.Range("E35").Value = '  -1.39%  '
$ws.Range("E36").Value = '  -1.03%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.916'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -3.00%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0164'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -0.19%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.526'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +0.87%  '
$ws.Range("E40").Value = '  -0.35%  '
$ws.Range("E41").Value = '  -0.14%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.995'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.73%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.38'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +3.35%  '
$ws.Range("E44").Value = '  +0.05%  '
$ws.Range("B45").Value = 'RenderToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.76'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -2.27%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '63.16'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.32%  '
$ws.Range("D47").Value = '1.698.79'
$ws.Range("E47").Value = '  +0.10%  '
$ws.Range("E48").Value = '  +1.03%  '
$ws.Range("E49").Value = '  +2.86%  '
$ws.Range("D50").Value = '0.0₇0981'
$ws.Range("E50").Value = '  -1.51%  '
$ws.Range("E51").Value = '  +0.71%  '
